$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ValidLoginTest")

# Insert a new row at row 2 (pushes the existing "receptionist"/"physician"
# rows down to rows 3/4) and fill it with a new admin/pass test case.
$ws.Rows.Item(2).Insert() | Out-Null
$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "pass"
$ws.Range("C2").Value = "English (Indian)"
$ws.Range("D2").Value = "OpenEMR"

# Mirror the author's final selection (column D across the now 4 data rows).
$ws.Range("D2:D4").Select() | Out-Null
